$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B (duplicate "จังหวัด" column) - shifts C:L left to B:K
$ws.Columns("B").Delete()

$ws.Range("A2").Select()
